$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.852.36'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.599.20'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.08%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '209.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.18%  '
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.477'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.23%  '
$ws.Range('E8').Value = '  -2.65%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0610'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '17.82'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0786'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.64%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.821.30'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.599.21'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.04'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.508'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.855.40'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.54%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.61'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.85%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0712'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.44%  '
$ws.Range('E19').Value = '  +0.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '188.97'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.30'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.93'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.14%  '
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.59'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.36%  '
$ws.Range('E26').Value = '  -4.10%  '
$ws.Range('E27').Value = '  -3.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.50'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '14.92'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.91%  '
$ws.Range('E30').Value = '  -2.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0465'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.91%  '
$ws.Range('E32').Value = '  -2.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.55%  '
$ws.Range('E34').Value = '  -1.24%  '
$ws.Range('E35').Value = '  -2.68%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.104.18'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.75%  '
$ws.Range('E37').Value = '  -2.76%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.798'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.97%  '
$ws.Range('E39').Value = '  -2.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.494'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.67%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '95.56'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.733.86'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.06%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.06'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.740'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.12%  '
$ws.Range('E45').Value = '  -1.32%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '52.94'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.45'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.31%  '
$ws.Range('E48').Value = '  -3.32%  '
$ws.Range('E49').Value = '  -1.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.18%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.34'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.33%  '
